{"js": "// Part 3's closing paragraph currently ends the document and carries the\n// \"_GoBack\" bookmark Word leaves at the last edit position. We add one\n// final wrap-up sentence after it (separated by a blank paragraph, matching\n// the document's existing visual rhythm) and re-anchor \"_GoBack\" onto the\n// new last paragraph, since that is where Word would have left the cursor.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst lastParagraph = items[items.length - 1];\n\n// Drop the stale \"_GoBack\" bookmark before moving it.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Blank separator paragraph, then the new closing remark.\nconst blankParagraph = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nblankParagraph.insertParagraph(\n  \"With these changes, the agent seems to make it to its destination a bit more frequently than moving completely randomly.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Re-fetch the paragraph collection so the \"last paragraph\" anchor reflects\n// the freshly inserted content (a stale reference would re-collapse to the\n// start of the paragraph rather than after the new run).\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\nconst refreshedItems = refreshedParagraphs.items;\nconst finalParagraph = refreshedItems[refreshedItems.length - 1];\n\n// Re-create \"_GoBack\" collapsed at the end of the new final paragraph.\nconst endOfFinalParagraph = finalParagraph.getRange(\"End\");\nendOfFinalParagraph.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Part 3's closing paragraph ends the document and carries the \"_GoBack\"\n# bookmark Word leaves at the last edit position. We add one final wrap-up\n# sentence after it (separated by a blank paragraph, matching the\n# document's existing visual rhythm), inserting everything *before* the\n# existing \"_GoBack\" bookmark so it naturally stays anchored at the very\n# end of the document once we're done - exactly where Word would leave it.\n\n$d = $word.ActiveDocument\n\n# Blank separator paragraph, pushed in just before \"_GoBack\".\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Range.InsertBefore(\"`r\")\n\n# New closing remark paragraph break, again just before \"_GoBack\".\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Range.InsertBefore(\"`r\")\n\n# The new sentence itself, still inserted just before \"_GoBack\" so the\n# bookmark ends up collapsed right after it, same as before the edit.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Range.InsertBefore(\"With these changes, the agent seems to make it to its destination a bit more frequently than moving completely randomly.\")\n"}
